$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# ---------------------------------------------------------------------------
# 1. Remove the "Shill" (row 2) and "Recruit" (row 3) cards entirely. This
#    shifts every following row up by two.
# ---------------------------------------------------------------------------
$ws.Rows("2:3").Delete()

# ---------------------------------------------------------------------------
# 2. "Build Nukes" is now row 2: its gold cost becomes the literal text "X"
#    and its rules text is rewritten to reflect the variable cost.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "X"
$ws.Range("E2").Value = "Spend X uranium. Build X nukes"

# ---------------------------------------------------------------------------
# 3. "Nuke" is now row 3: its rules text is rewritten with new terrorism /
#    nuke wording.
# ---------------------------------------------------------------------------
$ws.Range("E3").Value = "Senate votes on whether to activate this card or not. " + $nl + "If the vote passes, spend a nuke, then select a tile that is not in or adjacent to a Nation that has nuke(s). Remove all military units and workers inside the selected tile and the tiles adjacent to it. Nothing can be placed on those tiles for two turns. Decrease relations with all Nations by 1. If a Nation has a tile affected by "

# ---------------------------------------------------------------------------
# 4. Column E narrows (text got shorter overall), which in turn reflows
#    several row heights.
# ---------------------------------------------------------------------------
$ws.Columns("E").ColumnWidth = 99.75

# ---------------------------------------------------------------------------
# 5. Re-apply the row heights that changed because of the re-wrapped text /
#    narrower column. Rows that already match are left untouched.
# ---------------------------------------------------------------------------
$ws.Rows(3).RowHeight = 60
$ws.Rows(4).RowHeight = 30
$ws.Rows(5).RowHeight = 30
$ws.Rows(8).RowHeight = 45
$ws.Rows(10).RowHeight = 60
$ws.Rows(13).RowHeight = 105
$ws.Rows(14).RowHeight = 30

# ---------------------------------------------------------------------------
# 6. Update the view: scroll so column E is visible at the left edge and the
#    active selection sits on E7 (previously E10, before the two rows were
#    removed).
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("E7").Select()
